# Updates cryptocurrency price/volume data to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.344.06'
$ws.Range("E2").Value = '  +0.97%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.431.02'
$ws.Range("E3").Value = '  +1.78%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.42'
$ws.Range("E5").Value = '  +0.97%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.22'
$ws.Range("E6").Value = '  +6.85%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.432.33'
$ws.Range("E7").Value = '  +1.87%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.478'
$ws.Range("E9").Value = '  +2.29%  '

# Row 10
$ws.Range("E10").Value = '  +0.83%  '

# Row 11
$ws.Range("E11").Value = '  +3.45%  '

# Row 12
$ws.Range("E12").Value = '  +1.99%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.016.44'
$ws.Range("E13").Value = '  +1.78%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.93'
$ws.Range("E14").Value = '  +7.58%  '

# Row 15
$ws.Range("E15").Value = '  -0.56%  '

# Row 16
$ws.Range("E16").Value = '  +1.90%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.432.10'
$ws.Range("E17").Value = '  +1.80%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.456.51'
$ws.Range("E18").Value = '  +1.01%  '

# Row 19
$ws.Range("E19").Value = '  +8.54%  '

# Row 20
$ws.Range("E20").Value = '  +3.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.42'
$ws.Range("E21").Value = '  +2.44%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '395.42'
$ws.Range("E22").Value = '  +6.43%  '

# Row 23
$ws.Range("E23").Value = '  +3.27%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.31'
$ws.Range("E24").Value = '  +5.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.25%  '

# Row 26
$ws.Range("E26").Value = '  +0.39%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000122'
$ws.Range("E27").Value = '  -0.09%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.572.46'
$ws.Range("E28").Value = '  +1.86%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.179'
$ws.Range("E29").Value = '  +2.85%  '

# Row 30
$ws.Range("E30").Value = '  +3.85%  '

# Row 31
$ws.Range("E31").Value = '  -0.01%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.19'
$ws.Range("E32").Value = '  +2.25%  '

# Row 33
$ws.Range("E33").Value = '  -8.12%  '

# Row 34
$ws.Range("E34").Value = '  +2.46%  '

# Row 35
$ws.Range("E35").Value = '  -0.01%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.92'
$ws.Range("E36").Value = '  +2.77%  '

# Row 37
$ws.Range("E37").Value = '  +3.93%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.459.32'
$ws.Range("E38").Value = '  +2.01%  '

# Row 39
$ws.Range("E39").Value = '  +1.71%  '

# Row 40
$ws.Range("E40").Value = '  +0.40%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '167.53'
$ws.Range("E41").Value = '  +1.84%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0785'
$ws.Range("E42").Value = '  +3.41%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.70'
$ws.Range("E43").Value = '  +4.78%  '

# Row 44
$ws.Range("E44").Value = '  +3.74%  '

# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.04%  '

# Row 46
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.73'
$ws.Range("E46").Value = '  +0.06%  '

# Row 47
$ws.Range("E47").Value = '  +3.79%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.19'
$ws.Range("E48").Value = '  +0.83%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.594.20'
$ws.Range("E49").Value = '  +3.40%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.15'
$ws.Range("E50").Value = '  -0.32%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.92'
$ws.Range("E51").Value = '  +2.86%  '
